{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 1 (0-based index) is \"User: uytdsdfyuiop\" -> change the user's\n// message to \"hi\".\nconst userParagraph = paragraphs.items[1];\nconst oldUserText = \"uytdsdfyuiop\";\nconst userMatches = userParagraph.search(oldUserText, { matchCase: true });\nuserMatches.load(\"items\");\nawait context.sync();\nuserMatches.items[0].insertText(\"hi\", Word.InsertLocation.replace);\n\n// Paragraph 2 (0-based index) is the bot's first reply -> replace it with the\n// new, shorter greeting.\nconst botParagraph = paragraphs.items[2];\nconst oldBotText =\n  \"Oh, I see! *adjusts glasses* Well, hello there! *smiles* It's nice to meet you! How can I help you today? \\uD83D\\uDE0A Do you have an appointment with our doctor? Or do you have any questions or concerns? Please feel free to share! \\uD83E\\uDD17\";\nconst newBotText =\n  \"Hello there! *smiling* It's great to meet you! How can I help you today? Do you have an appointment with our doctor? Or do you have any questions or concerns about your health? Please feel free to share anything with me, and I'll do my best to assist you.\";\nconst botMatches = botParagraph.search(oldBotText, { matchCase: true });\nbotMatches.load(\"items\");\nawait context.sync();\nbotMatches.items[0].insertText(newBotText, Word.InsertLocation.replace);\nawait context.sync();\n\n// Remove the rest of the conversation (the follow-up User/Bot/User/Bot\n// exchanges) \u2014 they are no longer part of this page.\nfor (let i = 6; i >= 3; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 2: \"User: uytdsdfyuiop\" -> \"User: hi\"\n# Paragraph 3: \"Bot: <long greeting>\" -> \"Bot: <new shorter greeting>\"\n# Replace only the second run's text (after the bold \"User: \"/\"Bot: \" label)\n# by addressing the sub-range that starts right after the label, leaving the\n# label run and its bold formatting untouched.\n\n$p2 = $d.Paragraphs.Item(2)\n$p2Range = $p2.Range\n$labelLen = (\"User: \").Length\n$p2TextRange = $d.Range($p2Range.Start + $labelLen, $p2Range.End)\n$p2TextRange.Text = \"hi\"\n\n$p3 = $d.Paragraphs.Item(3)\n$p3Range = $p3.Range\n$labelLen2 = (\"Bot: \").Length\n$p3TextRange = $d.Range($p3Range.Start + $labelLen2, $p3Range.End)\n$p3TextRange.Text = \"Hello there! *smiling* It's great to meet you! How can I help you today? Do you have an appointment with our doctor? Or do you have any questions or concerns about your health? Please feel free to share anything with me, and I'll do my best to assist you.\"\n\n# Delete the remaining paragraphs (the follow-up User/Bot/User/Bot exchanges)\n# that are no longer part of the conversation (now paragraphs 4 through 7).\n$pCount = $d.Paragraphs.Count\nif ($pCount -ge 7) {\n    $startPara = $d.Paragraphs.Item(4)\n    $endPara = $d.Paragraphs.Item(7)\n    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $range.Delete()\n}\n"}
